# PDN Analyzer - "Last updates and orrder"
# Reorders/updates the BOM table: drops the old C1,C2 / C3 bulk-cap rows,
# renumbers designators, adds a new U4 regulator + its caps, and adds new
# resistor rows (R17, R20) plus a split-out R5 row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Wipe the existing data rows (2:19) so we start from a clean slate.
#    (Rows 16:19 don't exist yet in the source sheet; clearing a range
#    that extends past the current used range is a no-op for those rows.)
# ---------------------------------------------------------------------
$ws.Range("A2:G19").ClearContents()

# ---------------------------------------------------------------------
# 2. Write the new cell values row by row.
#    NumberFormat "@" is set just before writing any text value that
#    Excel would otherwise auto-convert to a number (package codes like
#    "0402"/"0603" and the long numeric JLCPCB basic-part id) - the
#    formatting gets normalised again in step 3.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "C16, C17"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = "C5246508"
$ws.Range("E2").Value = "220µF 20mΩ"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "875105144008"

$ws.Range("A3").Value = "CN1, CN2"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "C434806"
$ws.Range("F3").Value = "73412-0110"
$ws.Range("G3").Value = "Molex"

$ws.Range("A4").Value = "J1"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "C5333437"
$ws.Range("F4").Value = "A2541HWR-2x6P"

$ws.Range("A5").Value = "U1"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "C22444686"
$ws.Range("F5").Value = "OPA891DR"
$ws.Range("G5").Value = "Texas Instruments"

$ws.Range("A6").Value = "U2, U3"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "C468516"
$ws.Range("F6").Value = "ADA4940-1ACPZ-R7"
$ws.Range("G6").Value = "Analog Devices"

$ws.Range("A7").Value = "U4"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "C2685819"
$ws.Range("F7").Value = "TPS7A3901DSCR"

$ws.Range("A8").Value = "C3, C13"
$ws.Range("B8").Value = 2
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0402"
$ws.Range("E8").Value = "0.5pF"
$ws.Range("F8").Value = "0.5pF"

$ws.Range("A9").Value = "C4, C5, C6, C7, C9, C10, C11, C12, C14, C15, C18, C19, C20, C21, C22, C23"
$ws.Range("B9").Value = 16
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0402"
$ws.Range("E9").Value = "2.2uF"
$ws.Range("F9").Value = "2.2uF"

$ws.Range("A10").Value = "C24, C25, C28, C30"
$ws.Range("B10").Value = 4
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0402"
$ws.Range("E10").Value = "22uF"
$ws.Range("F10").Value = "22uF"

$ws.Range("A11").Value = "C26, C27, C29"
$ws.Range("B11").Value = 3
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0402"
$ws.Range("E11").Value = "100nF"
$ws.Range("F11").Value = "100nF"

$ws.Range("A12").Value = "R1, R8, R9, R14, R18, R19"
$ws.Range("B12").Value = 6
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0402"
$ws.Range("E12").Value = "10kΩ"
$ws.Range("F12").Value = "10kΩ"

$ws.Range("A13").Value = "R3, R7, R10, R13"
$ws.Range("B13").Value = 4
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0402"
$ws.Range("E13").Value = "1kΩ"
$ws.Range("F13").Value = "1kΩ"

$ws.Range("A14").Value = "R4, R6, R11, R12"
$ws.Range("B14").Value = 4
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0402"
$ws.Range("E14").Value = "100Ω"
$ws.Range("F14").Value = "100Ω"

$ws.Range("A15").Value = "R17"
$ws.Range("B15").Value = 1
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0402"
$ws.Range("E15").Value = "11kΩ"
$ws.Range("F15").Value = "11kΩ"

$ws.Range("A16").Value = "R20"
$ws.Range("B16").Value = 1
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0402"
$ws.Range("E16").Value = "21kΩ"
$ws.Range("F16").Value = "21kΩ"

$ws.Range("A17").Value = "R2"
$ws.Range("B17").Value = 1
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0603"
$ws.Range("E17").Value = "0.5Ω"
$ws.Range("F17").Value = "0.5Ω"

$ws.Range("A18").Value = "R15, R16"
$ws.Range("B18").Value = 2
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0402"
$ws.Range("E18").Value = "0Ω"
$ws.Range("F18").Value = "0Ω"

$ws.Range("A19").Value = "R5"
$ws.Range("B19").Value = 1
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0603"
$ws.Range("E19").Value = "0Ω"
$ws.Range("F19").Value = "0Ω"

# ---------------------------------------------------------------------
# 3. Normalise formatting.
#    Every data cell in the table uses one of two borders-only styles:
#      - the "plain" style (used for numbers & blanks)
#      - the "quote-prefix" style (used for text entries)
#    Paste the plain style across the whole block first, then re-apply
#    the quote-prefix style only to the cells that actually hold text.
# ---------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("A2:G19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A2").Copy()
$textCellRanges = @("A2:A19", "C2:C6", "D8:D19", "E2", "E8:E19", "F2:F19", "G3", "G5:G6")
foreach ($rng in $textCellRanges) {
    $ws.Range($rng).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# C7 is a pasted-in JLCPCB id that (per the source file) kept the plain
# (non quote-prefixed) style even though it's text - restore that.
$ws.Range("B7").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
# ...but PasteSpecial(formats) also carries the numeric flag off B7; put
# the value back as text afterwards.
$ws.Range("C7").Value = "C2685819"

# ---------------------------------------------------------------------
# 4. Selection, to match the saved workbook state.
# ---------------------------------------------------------------------
$ws.Range("C22").Select()
